# Upload of the last version of the burndown chart
# Fill in the newly-tracked "Added Sp" (G) and "Burned Sp" (F) figures for
# Sprint 1 / Sprint 2, and derive the "Sp ideali" (E) value for Sprint 3
# from the running burndown formula already used in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Sprint 1 row (row 3): Added Sp = 0
$ws.Range("G3").Value = 0

# Sprint 2 row (row 4): Burned Sp = 76, Added Sp = 0
$ws.Range("F4").Value = 76
$ws.Range("G4").Value = 0

# Sprint 3 row (row 5): Sp ideali follows the same pattern as row 4 (E4=E3-F3)
$ws.Range("E5").Formula = "=E4-F4"

$ws.Range("T17").Select()
